$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial for every data row (2..234).
# The automatic daily refresh bumps every one of these values by +1 day
# (46074 -> 46075), while everything else in the sheet stays the same.
$lastRow = $ws.Cells.Item(1, 3).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v + 1
    }
}
